# Generate Report for Handoff
# Update the "b.md" row across the Overview / zh-cn / de-de sheets to
# reflect a new handoff: status moves from "Handed back: in sync with
# en-US" to "Ready for handoff", a new handoff file (b.632...xlf) is
# produced per locale, and the handoff timestamps are refreshed.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"

# ---------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de | Latest Handoff Date
# Row 3 is the "b.md" entry.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusReady
$wsOverview.Range("C3").Value = $statusReady
$wsOverview.Range("D3").Value = "2016-03-21 08:29:54"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 ("b.md") gets a new handoff file + datetime.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusReady
$wsZhCn.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("E3").Value = "2016-03-21 08:29:50"

foreach ($h in $wsZhCn.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
    }
}

# ---------------------------------------------------------------------
# de-de sheet: row 3 ("b.md") gets a new handoff file + datetime.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusReady
$wsDeDe.Range("D3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("E3").Value = "2016-03-21 08:29:54"

foreach ($h in $wsDeDe.Hyperlinks) {
    if ($h.Range.Address() -eq '$D$3') {
        $h.TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
    }
}
